# Updates league database rows: swap the full row contents (columns B..AD)
# between rows that belong to the same re-ordered group, while leaving the
# row-index column (A) untouched. This mirrors the commit:
# "Atualizacao de bases das ligas, do dia: 15-06-2024 as 21:10"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each inner array is a cycle of worksheet rows. The content (columns B..AD)
# of row[0] moves to row[-1], row[1] moves to row[0], etc. In other words,
# every row ends up holding the data that used to live in the *next* row of
# the cycle (wrapping around).
$cycles = @(
    @(35, 36),
    @(59, 60),
    @(65, 66),
    @(119, 120),
    @(121, 122, 123),
    @(183, 184),
    @(200, 201),
    @(215, 216, 217),
    @(226, 227),
    @(230, 231),
    @(252, 253),
    @(271, 272),
    @(293, 294, 296)
)

$firstCol = 2   # column B
$lastCol  = 30  # column AD

foreach ($cycle in $cycles) {
    $count = $cycle.Length

    # Snapshot the current B..AD values for every row in this cycle first,
    # since we are about to overwrite them.
    $snapshots = @()
    foreach ($r in $cycle) {
        $rng = $ws.Range($ws.Cells.Item($r, $firstCol), $ws.Cells.Item($r, $lastCol))
        $snapshots += , $rng.Value2
    }

    # Row i receives the snapshot that belonged to row (i+1), wrapping
    # around so the last row receives the first row's original data.
    for ($i = 0; $i -lt $count; $i++) {
        $destRow = $cycle[$i]
        $srcSnapshot = $snapshots[($i + 1) % $count]
        $rng = $ws.Range($ws.Cells.Item($destRow, $firstCol), $ws.Cells.Item($destRow, $lastCol))
        $rng.Value2 = $srcSnapshot
    }
}
